# Update calibration data with new costs
# Rows 100-107 and 114-115, columns J:AS (cols 10-45) all get the same
# constant value replaced with a new constant value for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    100 = 103148.5282
    101 = 177833.2718
    102 = 833078.8325
    103 = 29151.07325
    104 = 79345.02172999999
    105 = 159707.5197
    106 = 91638.67157999999
    107 = 1107453.899
    114 = 163.8802033
    115 = 13243461.46
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $range = $ws.Range("J$row`:AS$row")
    $range.Value = $value
}
